$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in row 1
$ws.Range("A1").Value = "Sample Name"
$ws.Range("B1").Value = "bar"

# Update values in column A for rows 4 and 5
$ws.Range("A4").Value = 5
$ws.Range("A5").Value = 6

# Update row 1 height (auto-fit due to rotated/longer text)
$ws.Rows.Item(1).RowHeight = 69

# Update selection to D6
$ws.Range("D6").Select()
